$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Add the new "IMC" data row (row 5)
# ---------------------------------------------------------------------------
$ws.Range("A5").Value2 = "IMC"
$ws.Range("B5").Value2 = "https://www.imc.com/us/search-careers?jobTypes=Intern&page=1"
$ws.Range("C5").Value2 = "https://www.imc.com/us/"
$ws.Range("D5").Value2 = "_13fp8yk6c dzswju0 dzswjuc ryl3ea0 ajj7g52 _1e5s7rk1"
$ws.Range("E5").Value2 = "_13fp8yk6c dzswju0 dzswju5"

# ---------------------------------------------------------------------------
# 2. Apply word-wrap formatting to the existing/new data rows
#    (D3:E3 already wrap - leave them untouched)
# ---------------------------------------------------------------------------
$ws.Range("A2").WrapText = $true
$ws.Range("D2:E2").WrapText = $true
$ws.Range("B2:C2").WrapText = $true

$ws.Range("A3:C3").WrapText = $true

$ws.Range("A4:E4").WrapText = $true
$ws.Range("A5:E5").WrapText = $true

# ---------------------------------------------------------------------------
# 3. Create the extra empty, but formatted, rows (6-55)
# ---------------------------------------------------------------------------
for ($r = 6; $r -le 55; $r++) {
    $ws.Range("A" + $r + ":E" + $r).WrapText = $true
}

# ---------------------------------------------------------------------------
# 4. Row heights for the rows whose content now wraps onto multiple lines
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 29.5
$ws.Rows.Item(5).RowHeight = 44.25

# ---------------------------------------------------------------------------
# 5. Sheet view / selection
# ---------------------------------------------------------------------------
$null = $ws.Range("B5").Select()

Write-Host "done"
